$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 276 (shifts existing rows 276-296 down to 277-297)
$ws.Rows(276).Insert()

# Fill the new row 276 with the new price-report record
$ws.Range("A276").Value = 1
$ws.Range("B276").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C276").Value = "Arica y Parinacota"
$ws.Range("D276").Value = 44753
$ws.Range("E276").Value = 15
$ws.Range("F276").Value = 100114013
$ws.Range("G276").Value = "Zanahoria"
$ws.Range("H276").Value = "Sin especificar"
$ws.Range("I276").Value = "Primera"
$ws.Range("J276").Value = 90
$ws.Range("K276").Value = 14000
$ws.Range("L276").Value = 15000
$ws.Range("M276").Value = 14500
$ws.Range("N276").Value = "`$/saco 25 kilos"
$ws.Range("O276").Value = "Valle de Camiña"
$ws.Range("P276").Value = 580
$ws.Range("Q276").Value = 25
$ws.Range("R276").Value = "Hortaliza"
